$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2020 data column (N), matching the formatting already used for
# the other year columns: N4 (header year) picks up the same style as the
# other plain year headers (K4), N5 (value) picks up the same style as the
# neighbouring year value (M5).
$xlPasteFormats = -4122

$ws.Range("K4").Copy()
$ws.Range("N4").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(4, 14).Value = 2020

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(5, 14).Value = 534

# Update the view: scroll so column E is the top-left visible column,
# and set the active selection to S10
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("S10").Select()
